$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.48%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.58%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.102"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.36%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07620"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.40%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.607"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.05%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'2.470"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.47%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9024"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.51%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1125"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'12.14%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1791"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.06%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09140"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.47%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04220"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-5.32%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1050"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.46%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001252"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.75%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005706"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.10%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.348"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.16%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.258"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.3293"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.79%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.658"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-5.80%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.19%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'1.43%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04064"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.84%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001246"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'0.004112"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.05%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D38").Value = "'0.02392"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'2.17%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05182"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.33%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007766"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.73%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-1.68%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007060"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'9.34%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'0.00%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007730"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-5.72%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3079"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-7.32%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'7.32%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05582"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1,518.68%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.05%"
$ws.Range("E51").Style = "Normal"
